$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("展览")
if (1 -eq 1) {
    # --- Insert the new row 15 (南昌·童话镇国乙&鸢only) ---
    $wsA.Range("A15").EntireRow.Insert()
    $wsA.Range("A14").Copy()
    $wsA.Range("A15").PasteSpecial(-4122)

    # --- Insert 2 new rows at 24:25 (景德镇.. / 江西..) ---
    $wsA.Range("A24:A25").EntireRow.Insert()
    $wsA.Range("A23").Copy()
    $wsA.Range("A24:A25").PasteSpecial(-4122)

    # Row 3: F -> 2991 (南昌·AP动漫游戏嘉年华)
    $wsA.Range("F3").Value = 2991

    # Row 7: F -> 1635 (南昌·CM01动漫游戏博览会)
    $wsA.Range("F7").Value = 1635

    # Row 9: F -> 80 (新余·文旅国漫嘉年华暨BM次元盛典)
    $wsA.Range("F9").Value = 80

    # Row 10: F -> 29 (景德镇·宅舞联萌·随舞动漫派对（免费活动))
    $wsA.Range("F10").Value = 29

    # Row 11: F -> 1342 (南昌·创造力动漫游戏嘉年华1.0)
    $wsA.Range("F11").Value = 1342

    # Row 13: F -> 473 (赣州·第三届半夏动漫展)
    $wsA.Range("F13").Value = 473

    # Row 15 (new): 南昌·童话镇国乙&鸢only
    $wsA.Range("B15").NumberFormat = "@"
    $wsA.Range("B15").Value = "2024-04-05"
    $wsA.Range("B15").Style = "Normal"
    $wsA.Range("C15").Value = "南昌·童话镇国乙&鸢only"
    $wsA.Range("D15").Value = "赣江北大道新力外滩9号 百嘉喜宴"
    $wsA.Range("E15").Value = "2024.04.05 10:00-04.05 21:00"
    $wsA.Range("F15").Value = 1
    $wsA.Range("G15").Value = 118
    $wsA.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=83012"
    $wsA.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202403/iyQuodAz1710834305273.jpeg"

    # Row 17: F -> 60 (萍乡·2024DDS国漫盛典)
    $wsA.Range("F17").Value = 60

    # Row 19: F -> 88 (南昌·第二届漫拥动漫嘉年华mini)
    $wsA.Range("F19").Value = 88

    # Row 21: F -> 3095 (南昌·New World国潮动漫博览会)
    $wsA.Range("F21").Value = 3095

    # Row 22: F -> 378 (九江·第三届ACD动漫游戏嘉年华)
    $wsA.Range("F22").Value = 378

    # Row 23: F -> 102 (吉安·COMIC LIFE次元假日04)
    $wsA.Range("F23").Value = 102

    # Row 24 (new): 景德镇·第十四届瓷都ACG动漫游戏博览会
    $wsA.Range("B24").NumberFormat = "@"
    $wsA.Range("B24").Value = "2024-05-01"
    $wsA.Range("B24").Style = "Normal"
    $wsA.Range("C24").Value = "景德镇·第十四届瓷都ACG动漫游戏博览会"
    $wsA.Range("D24").Value = "新厂西路315号 陶溪川发布大厅"
    $wsA.Range("E24").Value = "2024.05.01 10:00-05.02 17:00"
    $wsA.Range("F24").Value = 1
    $wsA.Range("G24").Value = 50
    $wsA.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83016"
    $wsA.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202403/c0q8seJL1710835930052.png"

    # Row 25 (new): 江西·第二十二届九江ACJJ国际动漫展
    $wsA.Range("B25").NumberFormat = "@"
    $wsA.Range("B25").Value = "2024-05-01"
    $wsA.Range("B25").Style = "Normal"
    $wsA.Range("C25").Value = "江西·第二十二届九江ACJJ国际动漫展"
    $wsA.Range("D25").Value = "体育路九江市体育中心-体育馆 九江市体育中心"
    $wsA.Range("E25").Value = "2024.05.01 09:00-05.02 17:00"
    $wsA.Range("F25").Value = 2
    $wsA.Range("G25").Value = "不可售"
    $wsA.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=83004"
    $wsA.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202403/lFThDvkh1710829330909.jpeg"

    # Row 27: F -> 88 (南昌·代号鸢盛花行only)
    $wsA.Range("F27").Value = 88

    # --- Re-number column A (index) for all data rows: A(r) = r-1 ---
    for ($r = 2; $r -le 27; $r++) {
        $wsA.Cells.Item($r, 1).Value = $r - 1
    }
}

$wsB = $wb.Worksheets.Item("全部类型")
if (1 -eq 1) {
    # --- Insert the new row 15 (南昌·童话镇国乙&鸢only) ---
    $wsB.Range("A15").EntireRow.Insert()
    $wsB.Range("A14").Copy()
    $wsB.Range("A15").PasteSpecial(-4122)

    # --- Insert 2 new rows at 24:25 (景德镇.. / 江西..) ---
    $wsB.Range("A24:A25").EntireRow.Insert()
    $wsB.Range("A23").Copy()
    $wsB.Range("A24:A25").PasteSpecial(-4122)

    # Row 3: F -> 2991 (南昌·AP动漫游戏嘉年华)
    $wsB.Range("F3").Value = 2991

    # Row 7: F -> 1635 (南昌·CM01动漫游戏博览会)
    $wsB.Range("F7").Value = 1635

    # Row 9: F -> 80 (新余·文旅国漫嘉年华暨BM次元盛典)
    $wsB.Range("F9").Value = 80

    # Row 10: F -> 29 (景德镇·宅舞联萌·随舞动漫派对（免费活动))
    $wsB.Range("F10").Value = 29

    # Row 11: F -> 1342 (南昌·创造力动漫游戏嘉年华1.0)
    $wsB.Range("F11").Value = 1342

    # Row 13: F -> 473 (赣州·第三届半夏动漫展)
    $wsB.Range("F13").Value = 473

    # Row 15 (new): 南昌·童话镇国乙&鸢only
    $wsB.Range("B15").NumberFormat = "@"
    $wsB.Range("B15").Value = "2024-04-05"
    $wsB.Range("B15").Style = "Normal"
    $wsB.Range("C15").Value = "南昌·童话镇国乙&鸢only"
    $wsB.Range("D15").Value = "赣江北大道新力外滩9号 百嘉喜宴"
    $wsB.Range("E15").Value = "2024.04.05 10:00-04.05 21:00"
    $wsB.Range("F15").Value = 1
    $wsB.Range("G15").Value = 118
    $wsB.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=83012"
    $wsB.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202403/iyQuodAz1710834305273.jpeg"

    # Row 17: F -> 60 (萍乡·2024DDS国漫盛典)
    $wsB.Range("F17").Value = 60

    # Row 19: F -> 88 (南昌·第二届漫拥动漫嘉年华mini)
    $wsB.Range("F19").Value = 88

    # Row 21: F -> 3095 (南昌·New World国潮动漫博览会)
    $wsB.Range("F21").Value = 3095

    # Row 22: F -> 378 (九江·第三届ACD动漫游戏嘉年华)
    $wsB.Range("F22").Value = 378

    # Row 23: F -> 102 (吉安·COMIC LIFE次元假日04)
    $wsB.Range("F23").Value = 102

    # Row 24 (new): 景德镇·第十四届瓷都ACG动漫游戏博览会
    $wsB.Range("B24").NumberFormat = "@"
    $wsB.Range("B24").Value = "2024-05-01"
    $wsB.Range("B24").Style = "Normal"
    $wsB.Range("C24").Value = "景德镇·第十四届瓷都ACG动漫游戏博览会"
    $wsB.Range("D24").Value = "新厂西路315号 陶溪川发布大厅"
    $wsB.Range("E24").Value = "2024.05.01 10:00-05.02 17:00"
    $wsB.Range("F24").Value = 1
    $wsB.Range("G24").Value = 50
    $wsB.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83016"
    $wsB.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202403/c0q8seJL1710835930052.png"

    # Row 25 (new): 江西·第二十二届九江ACJJ国际动漫展
    $wsB.Range("B25").NumberFormat = "@"
    $wsB.Range("B25").Value = "2024-05-01"
    $wsB.Range("B25").Style = "Normal"
    $wsB.Range("C25").Value = "江西·第二十二届九江ACJJ国际动漫展"
    $wsB.Range("D25").Value = "体育路九江市体育中心-体育馆 九江市体育中心"
    $wsB.Range("E25").Value = "2024.05.01 09:00-05.02 17:00"
    $wsB.Range("F25").Value = 2
    $wsB.Range("G25").Value = "不可售"
    $wsB.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=83004"
    $wsB.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202403/lFThDvkh1710829330909.jpeg"

    # Row 27: F -> 88 (南昌·代号鸢盛花行only)
    $wsB.Range("F27").Value = 88

    # --- Re-number column A (index) for all data rows: A(r) = r-1 ---
    for ($r = 2; $r -le 27; $r++) {
        $wsB.Cells.Item($r, 1).Value = $r - 1
    }
}

Write-Host "Edit complete."